$d = $word.ActiveDocument

# Locate the paragraph that starts the block of paragraphs to remove:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   (empty paragraph)
#   (empty paragraph, page-break-before)
#   (empty paragraph)
# These four paragraphs, sitting right after the "Requisitos" /
# "LOB1036: ..." paragraphs, are deleted in their entirety while the
# paragraph mark that follows them (the final, empty page-break-before
# paragraph that precedes the section properties) is left untouched.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $targetStart = $searchRange.Start

    # Resolve the paragraph index (1-based) that contains the match.
    $startIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Start -le $targetStart -and $candidate.Range.End -gt $targetStart) {
            $startIndex = $i
            break
        }
    }

    if ($startIndex -gt 0) {
        $paragraphsToRemove = 4
        $endIndex = $startIndex + $paragraphsToRemove - 1

        $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
        $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End

        $delRange = $d.Range($rangeStart, $rangeEnd)
        $delRange.Delete()
    }
}
